$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet view selection now spans the full column A (test ranges were renamed/refreshed)
$ws.Range("A1:A1048576").Select()

# Refresh the r2 score matrix (C2:M18). Each inner array is [rowNumber, C..M values]
$r2Rows = @(
    @(2, -0.09, 0.03, 0.01, 0.05, -0.05, 0, 0.01, 0.05, -0.02, 0.04, -0.03),
    @(3, -0.09, 0.03, 0, 0.05, -0.05, 0, 0.02, 0.06, -0.02, 0.05, 0.07),
    @(4, -0.09, 0.03, -0.01, 0.05, -0.05, 0, 0.02, 0.06, -0.02, 0.05, 0.07),
    @(5, -0.09, 0.03, -0.01, 0.05, -0.05, 0, 0.02, 0.06, -0.02, 0.05, 0.08),
    @(6, -0.13, -0.1, -0.3, -0.06, -0.05, -0.04, 0.02, -0.31, 0.09, -0.05, -0.03),
    @(7, -0.13, -0.14, -0.3, -0.98, -18.53, -0.04, -0.04, -0.31, 0.08, -0.05, -0.05),
    @(8, -0.13, -0.15, -0.3, -0.12, -27.1, -0.04, 0.02, -0.31, 0.09, -0.05, -0.05),
    @(9, -0.13, -0.17, -0.3, -0.14, -0.32, -0.04, -0.02, -0.31, 0.09, -0.05, -0.05),
    @(10, -0.1, 0.04, 0.01, 0.08, 0.07, -0.01, 0.05, 0.06, -0.04, 0.05, 0.11),
    @(11, -0.1, 0.04, 0, 0.08, 0.09, -0.01, 0.05, 0.07, -0.04, 0.06, 0.12),
    @(12, -0.1, 0.04, -0.01, 0.08, 0.09, -0.01, 0.05, 0.07, -0.04, 0.06, 0.12),
    @(13, -0.1, 0.04, -0.01, 0.07, 0.09, -0.01, 0.05, 0.07, -0.04, 0.06, 0.12),
    @(14, -0.14, -0.11, -0.31, -0.15, -0.06, -0.06, 0.02, -0.31, 0.07, -0.06, -0.02),
    @(15, -0.14, -0.15, -0.31, -0.15, -31.38, -0.06, -0.03, -0.31, 0.07, -0.06, -0.03),
    @(16, -0.14, -0.16, -0.31, -3.54, -50.05, -0.06, 0.01, -0.31, 0.07, -0.06, -0.06),
    @(17, -0.14, -0.18, -0.31, -0.15, -1.33, -0.06, -0.03, -0.31, 0.07, -0.06, -0.06),
    @(18, -0.09, 0.04, 0.01, 0.08, 0.09, 0, 0.05, 0.07, 0.09, 0.06, 0.12)
)

foreach ($entry in $r2Rows) {
    $rowNum = $entry[0]
    for ($i = 1; $i -lt $entry.Length; $i++) {
        $ws.Cells.Item($rowNum, 2 + $i).Value = $entry[$i]
    }
}
